$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for week ending 2021-10-24 (rows 1706-1735)
# Column A uses the same date display style as the rest of the sheet (yyyy-mm-dd hh:mm:ss).
$ws.Range("A1706:A1735").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$rows = @(
    @{A=44487; B="'4182031"; C=3011; D="Order 4182031 Card(Stripe)"; E=$null; F=685.71},
    @{A=44487; B="'4182031"; C=2611; D="Order 4182031 Card(Stripe)"; E=$null; F=82.29},
    @{A=44487; B="'4182031"; C=1930; D="Order 4182031 Card(Stripe)"; E=768; F=$null},
    @{A=44487; B=$null; C=4010; D="WILLYS RISSNE K6885"; E=232.23; F=$null},
    @{A=44487; B=$null; C=2645; D="WILLYS RISSNE K6885"; E=27.87; F=$null},
    @{A=44487; B=$null; C=1930; D="WILLYS RISSNE K6885"; E=$null; F=260.1},
    @{A=44488; B="'7190943"; C=3011; D="Order 7190943 Card(Stripe)"; E=$null; F=650.89},
    @{A=44488; B="'7190943"; C=2611; D="Order 7190943 Card(Stripe)"; E=$null; F=78.11},
    @{A=44488; B="'7190943"; C=1930; D="Order 7190943 Card(Stripe)"; E=729; F=$null},
    @{A=44489; B=$null; C=6570; D="Pris banktjänster enligt faktura"; E=64.5; F=$null},
    @{A=44489; B=$null; C=$null; D="Pris banktjänster enligt faktura"; E=0; F=$null},
    @{A=44489; B=$null; C=1930; D="Pris banktjänster enligt faktura"; E=$null; F=64.5},
    @{A=44490; B="'3211947"; C=3011; D="Order 3211947 Card(Stripe)"; E=$null; F=423.21},
    @{A=44490; B="'3211947"; C=2611; D="Order 3211947 Card(Stripe)"; E=$null; F=50.79},
    @{A=44490; B="'3211947"; C=1930; D="Order 3211947 Card(Stripe)"; E=474; F=$null},
    @{A=44490; B=$null; C=4010; D="NGROCERIES K0135"; E=625.66; F=$null},
    @{A=44490; B=$null; C=2645; D="NGROCERIES K0135"; E=75.08; F=$null},
    @{A=44490; B=$null; C=1930; D="NGROCERIES K0135"; E=$null; F=700.74},
    @{A=44490; B=$null; C=4010; D="M&S RB BROMMA K0135"; E=326.24; F=$null},
    @{A=44490; B=$null; C=2645; D="M&S RB BROMMA K0135"; E=39.15; F=$null},
    @{A=44490; B=$null; C=1930; D="M&S RB BROMMA K0135"; E=$null; F=365.39},
    @{A=44491; B="'5221941"; C=3011; D="Order 5221941 Card(Stripe)"; E=$null; F=1062.5},
    @{A=44491; B="'5221941"; C=2611; D="Order 5221941 Card(Stripe)"; E=$null; F=127.5},
    @{A=44491; B="'5221941"; C=1930; D="Order 5221941 Card(Stripe)"; E=1190; F=$null},
    @{A=44492; B="'6231033"; C=3011; D="Order 6231033 Swish +46769332411"; E=$null; F=423.21},
    @{A=44492; B="'6231033"; C=2611; D="Order 6231033 Swish +46769332411"; E=$null; F=50.79},
    @{A=44492; B="'6231033"; C=1930; D="Order 6231033 Swish +46769332411"; E=474; F=$null},
    @{A=44491; B=$null; C=6400; D="FACEBK LEWEZ6XY62 K6885"; E=430; F=$null},
    @{A=44491; B=$null; C=$null; D="FACEBK LEWEZ6XY62 K6885"; E=0; F=$null},
    @{A=44491; B=$null; C=1930; D="FACEBK LEWEZ6XY62 K6885"; E=$null; F=430},
)

$startRow = 1706
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}
